# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account-statement (Estado de Cuenta) data for rows 16..47
# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @(16, "CC", "7919786",    "JORGE ANDRES GIRALDO SUAREZ", "2009", 35120, 878000),
    @(17, "CC", "1040510476", "EVA SANDRITH BARON LUCAS",     "2009", 35120, 878000),
    @(18, "CC", "73184471",   "DIOGENES ALFONSO PIÑA GUERRERO","2008", 35120, 878000),
    @(19, "CC", "45546890",   "KEILA ROSA GAVIRIA MARTINEZ",  "2008", 35120, 878000),
    @(20, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2105", 28096, 878000),
    @(21, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2104", 35120, 878000),
    @(22, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2103", 35120, 878000),
    @(23, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2102", 35120, 878000),
    @(24, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2101", 35120, 878000),
    @(25, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2012", 35120, 878000),
    @(26, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2011", 35120, 878000),
    @(27, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2010", 35120, 878000),
    @(28, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2009", 35120, 878000),
    @(29, "CC", "1049393176", "ALEXANDER CARRILLO VALENCIA",  "2008", 7024,  878000),
    @(30, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2105", 28096, 878000),
    @(31, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2104", 35120, 878000),
    @(32, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2103", 35120, 878000),
    @(33, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2102", 35120, 878000),
    @(34, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2101", 35120, 878000),
    @(35, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2012", 35120, 878000),
    @(36, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2011", 35120, 878000),
    @(37, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2010", 35120, 878000),
    @(38, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2009", 35120, 878000),
    @(39, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2008", 35120, 878000),
    @(40, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2007", 35120, 878000),
    @(41, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2006", 35120, 878000),
    @(42, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2005", 35120, 878000),
    @(43, "CC", "1102883290", "ANUAR ANDRES CASTELLAR RIOS",  "2004", 35120, 878000),
    @(44, "CC", "13197398",   "OMAR RINCON URBINA",           "2009", 35120, 878000),
    @(45, "CC", "13197398",   "OMAR RINCON URBINA",           "2008", 35120, 878000),
    @(46, "CC", "13197398",   "OMAR RINCON URBINA",           "2007", 35120, 878000),
    @(47, "CC", "13197398",   "OMAR RINCON URBINA",           "2006", 14048, 878000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
